$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginUser")

# Row 3: the old "LoginWithPassword" test case is renamed to "LoginWithIncorrectEmail"
# (its Email/Password values are left as-is).
$ws.Range("A3").Value = "LoginWithIncorrectEmail"

# New row 4: "LoginWithoutPassword" test case (Email filled in, Password left blank).
$ws.Range("A4").Value = "LoginWithoutPassword"
$ws.Range("A4").HorizontalAlignment = -4108  # xlCenter, matching the rest of the table
$ws.Range("B4").Value = "test@gmail.com"

# Widen column B slightly to fit the new e-mail address.
$ws.Columns.Item(2).ColumnWidth = 15.66

# Match the author's final selection/cursor position.
$ws.Range("B4").Select()

# Page setup used when the sheet was last saved.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
